# Generate Report for handoff
# Adds a new "Ready for handoff" row (for file
# 19610356-b2fd-4781-a25d-de0dee7cda64.md, producing handoff artifact
# ffff594e8794-256f-480c-acc1-8f630d18b6d6.md) above the existing
# ".localization-config" row on every worksheet, and renames the old
# "a7121821-1a0d-4dfc-9a39-c640b6860ad5.md" handoff to the new guid.

$wb = $excel.ActiveWorkbook

$oldFile = "a7121821-1a0d-4dfc-9a39-c640b6860ad5.md"
$newFile = "19610356-b2fd-4781-a25d-de0dee7cda64.md"
$handoffFile = "ffff594e8794-256f-480c-acc1-8f630d18b6d6.md"
$status = "Ready for handoff"

$oldFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8daa92ec74ac0976e3165cf36a85b0787a765527/e2e/a7121821-1a0d-4dfc-9a39-c640b6860ad5.md"
$newFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8daa92ec74ac0976e3165cf36a85b0787a765527/e2e/19610356-b2fd-4781-a25d-de0dee7cda64.md"
$handoffFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8daa92ec74ac0976e3165cf36a85b0787a765527/e2e/ffff594e8794-256f-480c-acc1-8f630d18b6d6.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/792a8bdbfe0fbaaaea3a301e2dc230eef379c6f3/.localization-config"

$zhXlf = "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf"
$deXlf = "19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a32d951201a97dda6ccd72c741106d657d07d3d4/localization/zh-cn/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a32d951201a97dda6ccd72c741106d657d07d3d4/localization/de-de/19610356-b2fd-4781-a25d-de0dee7cda64.a32d951201a97dda6ccd72c741106d657d07d3d4.de-de.xlf"

$zhHandoffDatetime = "2016-01-25 14:04:38"
$deHandoffDatetime = "2016-01-25 14:04:47"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

# Make room for the new "handoff" row by shifting row 3 (.localization-config)
# down to row 4.
$wsOverview.Rows.Item(3).Insert()

$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status

$wsOverview.Range("A3").Value = $handoffFile
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $newFileUrl, "", "", $newFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $handoffFileUrl, "", "", $handoffFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $configUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": same layout, different target-language data
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Xlf = $zhXlf; XlfUrl = $zhXlfUrl; Datetime = $zhHandoffDatetime },
    @{ Name = "de-de"; Xlf = $deXlf; XlfUrl = $deXlfUrl; Datetime = $deHandoffDatetime }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    $ws.Hyperlinks.Delete()

    # Shift row 3 (.localization-config) down to row 4.
    $ws.Rows.Item(3).Insert()

    # Row 2: rename the handed-off source file and mark it ready, with its
    # handoff target/file/datetime filled in.
    $ws.Range("A2").Value = $newFile
    $ws.Range("B2").Value = $status
    $ws.Range("C2").Value = $lang.Xlf
    $ws.Range("D2").Value = $lang.Datetime
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = "Include"

    # Row 3 (new): the handoff artifact file, also ready with the same
    # handoff target/file/datetime.
    $ws.Range("A3").Value = $handoffFile
    $ws.Range("B3").Value = $status
    $ws.Range("C3").Value = $lang.Xlf
    $ws.Range("D3").Value = $lang.Datetime
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = "Include"

    # Row 4 (was row 3): .localization-config, unchanged other than shifting.
    $ws.Range("A4").Value = ".localization-config"
    $ws.Range("B4").Value = "Not to be localized"
    $ws.Range("D4").Value = $epoch
    $ws.Range("G4").Value = $epoch
    $ws.Range("H4").Value = "Ignored"

    $ws.Hyperlinks.Add($ws.Range("A2"), $newFileUrl, "", "", $newFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $lang.XlfUrl, "", "", $lang.Xlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $handoffFileUrl, "", "", $handoffFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $lang.XlfUrl, "", "", $lang.Xlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", ".localization-config") | Out-Null
}
